# "Generate Report for Handoff"
# b.md was re-handed-off: update the Overview status plus the per-locale
# handoff file/datetime for the "b.md" row (row 3) on the zh-cn and de-de
# sheets, and repoint the matching hyperlink's displayed text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---- Overview sheet: row 3 (b.md) status in both locale columns ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---- zh-cn sheet: row 3 (b.md) new handoff file + datetime ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = $newStatus
$wsZh.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-01-29 02:22:39"

foreach ($h in $wsZh.Hyperlinks) {
    $rng = $h.Range()
    $addr = $rng.Address()
    if ($addr -eq '$C$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---- de-de sheet: row 3 (b.md) new handoff file + datetime ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = $newStatus
$wsDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("D3").Value = "2016-01-29 02:22:50"

foreach ($h in $wsDe.Hyperlinks) {
    $rng = $h.Range()
    $addr = $rng.Address()
    if ($addr -eq '$C$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
